$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells in columns B:E keep their literal string
# representation (Price/Volume columns hold "numbers" like "580.64" or
# "28.00" as plain text, e.g. European-style thousand-separated values like
# "61.877.93" that are not valid numerics). Force text format before writing.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.877.93"
$ws.Range("E2").Value = "  +1.04%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.461.57"
$ws.Range("E3").Value = "  +3.22%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.64"
$ws.Range("E5").Value = "  +1.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.35"
$ws.Range("E6").Value = "  +9.31%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.461.24"
$ws.Range("E7").Value = "  +3.23%  "

# Row 8
$ws.Range("E8").Value = "  +0.18%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  +1.26%  "

# Row 10
$ws.Range("E10").Value = "  +3.16%  "

# Row 11
$ws.Range("E11").Value = "  +2.03%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.389"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.048.43"
$ws.Range("E13").Value = "  +3.34%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.00"
$ws.Range("E14").Value = "  +8.23%  "

# Row 15
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.123"
$ws.Range("E15").Value = "  -0.14%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000175"
$ws.Range("E16").Value = "  +1.85%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.460.91"
$ws.Range("E17").Value = "  +3.60%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.966.98"
$ws.Range("E18").Value = "  +1.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("E19").Value = "  +8.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.36"
$ws.Range("E20").Value = "  +2.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.47"
$ws.Range("E21").Value = "  +2.44%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.26"
$ws.Range("E22").Value = "  +2.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.568"
$ws.Range("E23").Value = "  +3.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.590.59"
$ws.Range("E24").Value = "  +3.02%  "

# Row 26
$ws.Range("E26").Value = "  +0.97%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "72.51"
$ws.Range("E27").Value = "  +2.26%  "

# Row 28
$ws.Range("E28").Value = "  +0.35%  "

# Row 29
$ws.Range("E29").Value = "  +8.77%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.80"
$ws.Range("E30").Value = "  +4.37%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.54"
$ws.Range("E31").Value = "  -13.55%  "

# Row 32
$ws.Range("E32").Value = "  +0.02%  "

# Row 33
$ws.Range("E33").Value = "  +1.23%  "

# Row 34
$ws.Range("E34").Value = "  +2.40%  "

# Row 35
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.04"
$ws.Range("E36").Value = "  +2.21%  "

# Row 37
$ws.Range("E37").Value = "  +4.43%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.21"
$ws.Range("E38").Value = "  +0.05%  "

# Row 39
$ws.Range("E39").Value = "  +2.18%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.29"
$ws.Range("E40").Value = "  +0.98%  "

# Row 41
$ws.Range("E41").Value = "  +3.66%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.08"
$ws.Range("E42").Value = "  +10.02%  "

# Row 43
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.798"
$ws.Range("E43").Value = "  +3.90%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.73"
$ws.Range("E44").Value = "  +1.98%  "

# Row 45
$ws.Range("E45").Value = "  +0.02%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.36"
$ws.Range("E46").Value = "  +2.28%  "

# Row 47
$ws.Range("E47").Value = "  +2.37%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.17"
$ws.Range("E48").Value = "  -2.29%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.611.32"
$ws.Range("E49").Value = "  +11.32%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.50"
$ws.Range("E50").Value = "  +2.33%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.95"
$ws.Range("E51").Value = "  +2.35%  "
